$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 26318486
$ws.Range("I8").Value = 29414438
$ws.Range("J8").Value = 2902.5
$ws.Range("K8").Value = 88243314
$ws.Range("L8").Value = 8707.5
$ws.Range("M8").Value = -88243175
$ws.Range("N8").Value = -8985.5
# Row 19
$ws.Range("H19").Value = 590
$ws.Range("I19").Value = 635
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 635
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = -460
$ws.Range("N19").Value = -850
# Row 51
$ws.Range("H51").Value = 3000
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3968
# Row 69
$ws.Range("H69").Value = 250004660
$ws.Range("J69").Value = 250004660
$ws.Range("L69").Value = 750013980
$ws.Range("N69").Value = -750015728
# Row 72
$ws.Range("H72").Value = 250004660
$ws.Range("J72").Value = 250004660
$ws.Range("L72").Value = 2250041940
$ws.Range("N72").Value = -2250050676
# Row 112
$ws.Range("H112").Value = 1433.375
$ws.Range("J112").Value = 1411.1666
$ws.Range("L112").Value = 4233.4998
$ws.Range("N112").Value = -6449.4998
# Row 118
$ws.Range("H118").Value = 1725
$ws.Range("I118").Value = 1644.5
$ws.Range("J118").Value = 2208
$ws.Range("K118").Value = 4933.5
$ws.Range("L118").Value = 6624
$ws.Range("M118").Value = -3276.5
$ws.Range("N118").Value = -9938
# Row 129
$ws.Range("H129").Value = 1603.6
$ws.Range("J129").Value = 2218.2
$ws.Range("L129").Value = 6654.599999999999
$ws.Range("N129").Value = -16654.6
# Row 132
$ws.Range("H132").Value = 2223613.8
$ws.Range("I132").Value = 1423.2273
$ws.Range("K132").Value = 4269.6819
$ws.Range("M132").Value = -1739.6819
# Row 133
$ws.Range("H133").Value = 116101.43
$ws.Range("J133").Value = 116101.43
$ws.Range("L133").Value = 116101.43
$ws.Range("N133").Value = -126221.43
# Row 135
$ws.Range("H135").Value = 2262.875
$ws.Range("I135").Value = 1934
$ws.Range("K135").Value = 17406
$ws.Range("M135").Value = -14871
# Row 138
$ws.Range("H138").Value = 1920524.6
$ws.Range("I138").Value = 6131.6665
$ws.Range("J138").Value = 3007072
$ws.Range("K138").Value = 18394.9995
$ws.Range("L138").Value = 9021216
$ws.Range("M138").Value = -13254.9995
$ws.Range("N138").Value = -9031496
# Row 141
$ws.Range("H141").Value = 13798.429
$ws.Range("I141").Value = 13798.429
$ws.Range("K141").Value = 41395.287
$ws.Range("M141").Value = -36215.287

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10312.452
$ws.Range("I32").Value = 7431
$ws.Range("J32").Value = 27601.166
$ws.Range("K32").Value = 7431
$ws.Range("L32").Value = 27601.166
$ws.Range("M32").Value = -7144
$ws.Range("N32").Value = -28175.166
# Row 61
$ws.Range("H61").Value = 281553.03
$ws.Range("I61").Value = 2973.3096
$ws.Range("K61").Value = 2973.3096
$ws.Range("M61").Value = -2761.3096
# Row 130
$ws.Range("H130").Value = 44607.25
$ws.Range("J130").Value = 44607.25
$ws.Range("L130").Value = 44607.25
$ws.Range("N130").Value = -54647.25
# Row 132
$ws.Range("H132").Value = 1319.5122
$ws.Range("I132").Value = 1296.079
$ws.Range("K132").Value = 3888.237
$ws.Range("M132").Value = -1358.237
# Row 136
$ws.Range("H136").Value = 281553.03
$ws.Range("I136").Value = 2973.3096
$ws.Range("K136").Value = 8919.9288
$ws.Range("M136").Value = -6369.9288

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5792.974
$ws.Range("I31").Value = 1425
$ws.Range("K31").Value = 1425
$ws.Range("M31").Value = -1130
# Row 34
$ws.Range("H34").Value = 5792.974
$ws.Range("I34").Value = 1425
$ws.Range("K34").Value = 1425
$ws.Range("M34").Value = -1223
# Row 58
$ws.Range("H58").Value = 2874.476
$ws.Range("J58").Value = 3190.6667
$ws.Range("L58").Value = 3190.6667
$ws.Range("N58").Value = -3596.6667
# Row 80
$ws.Range("H80").Value = 50058
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -52246
# Row 83
$ws.Range("H83").Value = 50058
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -161232
# Row 122
$ws.Range("H122").Value = 47621184
$ws.Range("I122").Value = 3576.6667
$ws.Range("J122").Value = 83334390
$ws.Range("K122").Value = 10730.0001
$ws.Range("L122").Value = 250003170
$ws.Range("M122").Value = -8280.000100000001
$ws.Range("N122").Value = -250008070
# Row 132
$ws.Range("H132").Value = 2195.125
$ws.Range("I132").Value = 1524.25
$ws.Range("K132").Value = 4572.75
$ws.Range("M132").Value = -2042.75
# Row 136
$ws.Range("H136").Value = 2874.476
$ws.Range("J136").Value = 3190.6667
$ws.Range("L136").Value = 9572.000100000001
$ws.Range("N136").Value = -14672.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 119
$ws.Range("H119").Value = 1299
$ws.Range("I119").Value = 1299
$ws.Range("K119").Value = 3897
$ws.Range("M119").Value = 941
# Row 120
$ws.Range("H120").Value = 124879.22
$ws.Range("I120").Value = 139739.12
$ws.Range("J120").Value = 6000
$ws.Range("K120").Value = 419217.36
$ws.Range("L120").Value = 18000
$ws.Range("M120").Value = -414379.36
$ws.Range("N120").Value = -27676

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 742504500
$ws.Range("J15").Value = 742504500
$ws.Range("L15").Value = 742504500
$ws.Range("N15").Value = -742505076
# Row 55
$ws.Range("H55").Value = 10841.333
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 10841.333
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 10841.333
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -11495.333
# Row 62
$ws.Range("H62").Value = 600019200
$ws.Range("J62").Value = 750012500
$ws.Range("L62").Value = 750012500
$ws.Range("N62").Value = -750013872
# Row 65
$ws.Range("H65").Value = 600019200
$ws.Range("J65").Value = 750012500
$ws.Range("L65").Value = 2250037500
$ws.Range("N65").Value = -2250044364
# Row 81
$ws.Range("H81").Value = 742504500
$ws.Range("J81").Value = 742504500
$ws.Range("L81").Value = 742504500
$ws.Range("N81").Value = -742506496
# Row 84
$ws.Range("H84").Value = 742504500
$ws.Range("J84").Value = 742504500
$ws.Range("L84").Value = 2227513500
$ws.Range("N84").Value = -2227523484
# Row 102
$ws.Range("H102").Value = 1379.7059
$ws.Range("I102").Value = 1063.7
$ws.Range("K102").Value = 1063.7
$ws.Range("M102").Value = 558.3
# Row 107
$ws.Range("H107").Value = 2862.3635
$ws.Range("I107").Value = 996.75
$ws.Range("K107").Value = 996.75
$ws.Range("M107").Value = 923.25
# Row 132
$ws.Range("H132").Value = 25002192
$ws.Range("I132").Value = 25642760
$ws.Range("K132").Value = 76928280
$ws.Range("M132").Value = -76925750

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3516.9333
$ws.Range("I7").Value = 3472.2
$ws.Range("K7").Value = 3472.2
$ws.Range("M7").Value = -3360.2
# Row 16
$ws.Range("H16").Value = 706.3043
$ws.Range("J16").Value = 450
$ws.Range("L16").Value = 450
$ws.Range("N16").Value = -790
# Row 40
$ws.Range("H40").Value = 3896.3333
$ws.Range("I40").Value = 3134
$ws.Range("K40").Value = 3134
$ws.Range("M40").Value = -2998
# Row 126
$ws.Range("H126").Value = 3516.9333
$ws.Range("I126").Value = 3472.2
$ws.Range("K126").Value = 10416.6
$ws.Range("M126").Value = -7946.599999999999
# Row 137
$ws.Range("I137").Value = 93000
$ws.Range("J137").Value = 92000
$ws.Range("K137").Value = 93000
$ws.Range("L137").Value = 92000
$ws.Range("M137").Value = -87900
$ws.Range("N137").Value = -102200

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 6668310
$ws.Range("I126").Value = 1386.15
$ws.Range("J126").Value = 33336004
$ws.Range("K126").Value = 4158.450000000001
$ws.Range("L126").Value = 100008012
$ws.Range("M126").Value = -1688.450000000001
$ws.Range("N126").Value = -100012952
# Row 132
$ws.Range("H132").Value = 15153069
$ws.Range("J132").Value = 849.25
$ws.Range("L132").Value = 2547.75
$ws.Range("N132").Value = -7607.75
# Row 133
$ws.Range("H133").Value = 30547.555
$ws.Range("J133").Value = 30547.555
$ws.Range("L133").Value = 30547.555
$ws.Range("N133").Value = -40667.555

Write-Host "applied market data update"